$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.358.44"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").Value = "2.986.90"
$ws.Range("E3").Value = "  +2.41%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'384.57"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").Value = "'103.03"
$ws.Range("E6").Value = "  +0.35%  "

# Row 7
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
$ws.Range("D10").Value = "'36.78"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").Value = "'0.0843"
$ws.Range("E12").Value = "  +1.00%  "

# Row 13
$ws.Range("D13").Value = "3.449.08"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14
$ws.Range("D14").Value = "'18.19"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("E15").Value = "  +2.39%  "

# Row 16
$ws.Range("D16").Value = "2.984.19"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17
$ws.Range("D17").Value = "'0.995"
$ws.Range("E17").Value = "  +7.10%  "

# Row 18
$ws.Range("D18").Value = "51.296.83"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("D19").Value = "'3.26"
$ws.Range("E19").Value = "  -4.29%  "

# Row 20
$ws.Range("E20").Value = "  +0.66%  "

# Row 21
$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +1.49%  "

# Row 23
$ws.Range("D23").Value = "'68.86"
$ws.Range("E23").Value = "  +1.00%  "

# Row 24
$ws.Range("D24").Value = "'262.38"
$ws.Range("E24").Value = "  +0.80%  "

# Row 25
$ws.Range("E25").Value = "  +5.03%  "

# Row 26
$ws.Range("D26").Value = "'8.25"
$ws.Range("E26").Value = "  +15.44%  "

# Row 27
$ws.Range("D27").Value = "'7.52"
$ws.Range("E27").Value = "  +10.62%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.168"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("E29").Value = "  +12.98%  "

# Row 30
$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D30").Value = "'4.12"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("D32").Value = "'25.84"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("D33").Value = "'9.86"
$ws.Range("E33").Value = "  +0.92%  "

# Row 34
$ws.Range("D34").Value = "'34.48"
$ws.Range("E34").Value = "  +1.62%  "

# Row 35
$ws.Range("D35").Value = "'50.92"
$ws.Range("E35").Value = "  -0.33%  "

# Row 36
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  -2.90%  "

# Row 37
$ws.Range("D37").Value = "'0.0447"
$ws.Range("E37").Value = "  +6.69%  "

# Row 38
$ws.Range("E38").Value = "  +0.00%  "

# Row 39
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("D40").Value = "'17.08"
$ws.Range("E40").Value = "  +1.28%  "

# Row 41
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +1.70%  "

# Row 42
$ws.Range("E42").Value = "  +1.52%  "

# Row 43
$ws.Range("D43").Value = "'1.80"
$ws.Range("E43").Value = "  -0.41%  "

# Row 44
$ws.Range("D44").Value = "'122.35"
$ws.Range("E44").Value = "  -0.28%  "

# Row 45
$ws.Range("D45").Value = "'21.56"
$ws.Range("E45").Value = "  +0.42%  "

# Row 46
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  +0.35%  "

# Row 47
$ws.Range("D47").Value = "'0.277"
$ws.Range("E47").Value = "  +3.42%  "

# Row 48
$ws.Range("D48").Value = "'2.35"
$ws.Range("E48").Value = "  +1.71%  "

# Row 49
$ws.Range("D49").Value = "2.026.98"
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$ws.Range("D50").Value = "'3.25"
$ws.Range("E50").Value = "  +3.39%  "

# Row 51
$ws.Range("D51").Value = "'0.0332"
$ws.Range("E51").Value = "  +5.14%  "
